# Applies the "fixing stuff to mmol" commit: appends 18 new sample rows
# (rows 109-126) to Sheet1, mirroring the existing table's layout, styles
# and formulas, and updates the sheet's selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New row data, transcribed from the appended rows in the worksheet.
#    "kind" marks which columns are populated:
#      full   -> A..I populated, plus J/K formulas            (rows 109-114)
#      nohjk  -> A..G + I populated, no H (and so no J/K)      (rows 115-120)
#      short  -> A..E + I populated only (no F/G/H/J/K)        (rows 121-126)
# ---------------------------------------------------------------------
$rows = @(
    @{ r=109; A=45548; B=13;   C=1078; D=1.0677; E=45560; F=1.2854; G=45560; H=1.2844; kind='full' },
    @{ r=110; A=45548; B=5;    C=549;  D=1.122;  E=45560; F=1.3373; G=45560; H=1.3361; kind='full' },
    @{ r=111; A=45548; B=15;   C=552;  D=1.1053; E=45560; F=1.3227; G=45560; H=1.3211; kind='full' },
    @{ r=112; A=45548; B='6a'; C=563;  D=1.0768; E=45560; F=1.2931; G=45560; H=1.2917; kind='full' },
    @{ r=113; A=45548; B=7;    C=556;  D=1.0656; E=45560; F=1.2803; G=45560; H=1.2802; kind='full' },
    @{ r=114; A=45548; B=3;    C=564;  D=1.1361; E=45560; F=1.3529; G=45560; H=1.3519; kind='full' },

    @{ r=115; A=45548; B=6;    C=561;  D=2.3589; E=45560; F=1.3017; G=45565; kind='nohjk' },
    @{ r=116; A=45548; B='5a'; C=560;  D=2.3701; E=45560; F=1.3005; G=45565; kind='nohjk' },
    @{ r=117; A=45548; B=9;    C=564;  D=2.323;  E=45560; F=1.2912; G=45565; kind='nohjk' },
    @{ r=118; A=45553; B=6.2;  C=535;  D=2.7661; E=45560; F=1.3335; G=45565; kind='nohjk' },
    @{ r=119; A=45553; B=5.6;  C=541;  D=2.7868; E=45560; F=12797;  G=45565; kind='nohjk' },
    @{ r=120; A=45553; B=5.2;  C=545;  D=2.3309; E=45560; F=1.3213; G=45565; kind='nohjk' },

    @{ r=121; A=45553; B=9.4;  C=526;  D=2.5399; E=45565; kind='short' },
    @{ r=122; A=45553; B=5.1;  C=536;  D=2.7149; E=45565; kind='short' },
    @{ r=123; A=45553; B=5.4;  C=546;  D=2.7216; E=45565; kind='short' },
    @{ r=124; A=45553; B=3.1;  C=440;  D=2.3493; E=45565; kind='short' },
    @{ r=125; A=45553; B=9.3;  C=424;  D=2.6429; E=45565; kind='short' },
    @{ r=126; A=45553; B=5.3;  C=512;  D=2.694;  E=45565; kind='short' }
)

$firstNew = 109
$lastNew  = 126

# ---------------------------------------------------------------------
# 2. Pre-copy cell formatting (number formats) onto the new ranges so the
#    new cells pick up the same styles already used elsewhere in the
#    sheet (date style for A, date-time style for E, date style(s) for
#    G) without minting any new style entries.
# ---------------------------------------------------------------------

# Column A (date, rows 109-126): same style as the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A$firstNew`:A$lastNew").PasteSpecial(-4122)

# Column E (date-time, rows 109-120): same style as the rest of column E.
$ws.Range("E2").Copy()
$ws.Range("E$firstNew`:E120").PasteSpecial(-4122)

# Column E (date only, rows 121-126): matches the plain-date style used
# later in column G (e.g. G101).
$ws.Range("G101").Copy()
$ws.Range("E121:E126").PasteSpecial(-4122)

# Column G (date-time, rows 109-114): matches G100's style.
$ws.Range("G100").Copy()
$ws.Range("G109:G114").PasteSpecial(-4122)

# Column G (date only, rows 115-120): matches G101's style.
$ws.Range("G101").Copy()
$ws.Range("G115:G120").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Fill in the values for every new row.
# ---------------------------------------------------------------------
foreach ($row in $rows) {
    $r = $row.r

    $ws.Cells.Item($r, 1).Value2 = $row.A        # A - Sampled date
    $ws.Cells.Item($r, 2).Value2 = $row.B        # B - ID
    $ws.Cells.Item($r, 3).Value2 = $row.C        # C - Volume (mL)
    $ws.Cells.Item($r, 4).Value2 = $row.D        # D - Tray weight (g)

    if ($row.ContainsKey('E')) {
        $ws.Cells.Item($r, 5).Value2 = $row.E    # E - Desicator
    }

    if ($row.kind -eq 'full' -or $row.kind -eq 'nohjk') {
        $ws.Cells.Item($r, 6).Value2 = $row.F    # F - Desicator (g)
        $ws.Cells.Item($r, 7).Value2 = $row.G    # G - Furance
    }

    if ($row.kind -eq 'full') {
        $ws.Cells.Item($r, 8).Value2 = $row.H    # H - Furnace (g)
    }

    # I - Analyzer, always "SH" for the new rows.
    $ws.Cells.Item($r, 9).Value2 = "SH"

    if ($row.kind -eq 'full') {
        # J - POC (mg), K - mg/L, mirroring the existing formulas.
        $ws.Cells.Item($r, 10).Formula = "=(F$r-H$r)*1000"
        $ws.Cells.Item($r, 11).Formula = "=J$r/(C$r/1000)"
    }
}

# ---------------------------------------------------------------------
# 4. Update the sheet's selection so the view matches the state left
#    after entering the new data (last cell touched: I126).
# ---------------------------------------------------------------------
$ws.Range("I$lastNew").Select() | Out-Null
